# Apply the data update described by the commit "data modified in IPL folder":
#  - Insert two new columns ("ownTeam", "oppTeam") right before the existing
#    "batsman" column.
#  - Populate the new ownTeam/oppTeam columns for the pre-existing rows.
#  - Append four new match rows (5-8) with the full set of columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space + dagger used in "Dinesh Karthik (c)\u2020" (matches the
# text already present in the sheet).
$nbsp = [char]0x00A0
$dagger = [char]0x2020
$batsman = "Dinesh Karthik" + $nbsp + "(c)" + $dagger

# 1. Insert two blank columns before column D ("batsman") to make room for
#    "ownTeam" and "oppTeam". This shifts batsman..sr from D..I to F..K.
$ws.Range("D1:E1").EntireColumn.Insert()

# 2. New header cells for the inserted columns.
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# 3. Fill ownTeam / oppTeam for the pre-existing rows (2-4).
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Chennai Super Kings"

$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "Rajasthan Royals"

$ws.Range("D4").Value = "Kolkata Knight Riders"
$ws.Range("E4").Value = "Kings XI Punjab"

# 4. Append the new rows (5-8) with full data across all columns A-K.
$newRows = @(
    @{ A=" Abu Dhabi";  B=" September 26 2020"; C="KKR won by 7 wickets (with 12 balls remaining)"; D="Kolkata Knight Riders"; E="Sunrisers Hyderabad";         F=$batsman; G="0";  H="3";  I="0"; J="0"; K="0.00" },
    @{ A=" Sharjah";    B=" October 03 2020";   C="Capitals won by 18 runs";                         D="Kolkata Knight Riders"; E="Delhi Capitals";              F=$batsman; G="6";  H="8";  I="0"; J="0"; K="75.00" },
    @{ A=" Abu Dhabi";  B=" September 23 2020"; C="Mumbai won by 49 runs";                            D="Kolkata Knight Riders"; E="Mumbai Indians";              F=$batsman; G="30"; H="23"; I="5"; J="0"; K="130.43" },
    @{ A=" Sharjah";    B=" October 12 2020";   C="RCB won by 82 runs";                                D="Kolkata Knight Riders"; E="Royal Challengers Bangalore"; F=$batsman; G="1";  H="2";  I="0"; J="0"; K="50.00" }
)

$rowIndex = 5
foreach ($row in $newRows) {
    $ws.Range("A$rowIndex").Value = $row.A
    $ws.Range("B$rowIndex").Value = $row.B
    $ws.Range("C$rowIndex").Value = $row.C
    $ws.Range("D$rowIndex").Value = $row.D
    $ws.Range("E$rowIndex").Value = $row.E
    $ws.Range("F$rowIndex").Value = $row.F

    # Numeric-looking columns must stay text (matches existing "str" typed
    # cells for totalRuns/totalBalls/total4s/total6s/sr), so force the
    # number format to Text before assigning the value.
    $numRange = $ws.Range("G$($rowIndex):K$($rowIndex)")
    $numRange.NumberFormat = "@"
    $ws.Range("G$rowIndex").Value = $row.G
    $ws.Range("H$rowIndex").Value = $row.H
    $ws.Range("I$rowIndex").Value = $row.I
    $ws.Range("J$rowIndex").Value = $row.J
    $ws.Range("K$rowIndex").Value = $row.K

    $rowIndex++
}
